$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: relocate the trailing summary rows down to make room for the two
# new log entries (work bottom-up so we never overwrite a row before it has
# been relocated).
$ws.Range("C41:F41").Cut($ws.Range("C51"))
$ws.Range("D40:E40").Cut($ws.Range("D50"))
$ws.Range("D39:E39").Cut($ws.Range("D49"))

# Step 2: drop the now-vacated rows entirely (Cut left style-only remnants
# behind, and these rows must not exist in the final sheet).
$ws.Range("A39:F41").Clear()

# Step 3: Cut() collapses formulas down to static values, so restore them.
$ws.Range("E49").Formula = "=SUM(E2:E39)"
$ws.Range("E50").Formula = "=E49 / 60"

# Step 4: add the two new diary rows, reusing the formatting of the most
# recent data row (37) so the same cell styles are referenced instead of
# minting new ones.
$ws.Range("A37:F37").Copy()
$ws.Range("A38:F38").PasteSpecial(-4122)
$ws.Range("A37:F37").Copy()
$ws.Range("A39:F39").PasteSpecial(-4122)

# Step 5: populate values. New unique text is entered in the same order the
# original author typed it in, so new shared-string entries land on the
# same indices as the authoritative workbook.
$ws.Range("F38").Value = "Research on Linear Laddre Solutions, Algorithms and Diagrams"
$ws.Range("A39").Value = "21.10.2022"
$ws.Range("C38").Value = "Linear Ladder"
$ws.Range("F39").Value = "Research on Linear Laddre Solutions, Algorithms and Diagrams and R2R Ladders"

$ws.Range("A38").Value = "20.10.2022"
$ws.Range("B38").Value = 0.95833333333333337
$ws.Range("D38").Value = "Research"
$ws.Range("E38").Value = 210

$ws.Range("B39").Value = 0.5
$ws.Range("C39").Value = "Logarithmic Ladder"
$ws.Range("D39").Value = "Research"
$ws.Range("E39").Value = 190

# Step 6: reflect the new extent/selection the way Excel would after this
# edit (dimension grows to F51, and the last touched cell is E50).
$ws.Range("E50").Select() | Out-Null

Write-Host "Added Linear/Logarithmic Ladder log entries (rows 38-39); summary moved to rows 49-51."
